# "semana 29 de 2025" -- add week-29 column (AF) to the weekly IRA/UCI
# revision sheet, and backfill the previously-missing provider name for
# row 54 (CLINICA MEDICA TURIN SAS).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- header: AF1 = "29" (stored as text, like the other week-number
#     headers in row 1, not as a number) --------------------------------
$ws.Range("AF1").Value = "'29"

# --- week-29 counts for column AF, one entry per data row that carries
#     a value in this week ---------------------------------------------
$af29 = @{
  2  = 0
  3  = 0
  4  = 0
  5  = 0
  6  = 1
  7  = 0
  8  = 0
  9  = 0
  12 = 0
  14 = 0
  15 = 0
  17 = 0
  23 = 0
  24 = 0
  25 = 0
  26 = 0
  27 = 0
  28 = 1
  29 = 1
  30 = 2
  31 = 0
  32 = 0
  34 = 0
  35 = 2
  36 = 0
  37 = 0
  38 = 0
  39 = 0
  40 = 0
  41 = 0
  42 = 0
  43 = 0
  44 = 0
  45 = 0
  46 = 0
  47 = 0
  48 = 0
  49 = 0
  50 = 0
  52 = 0
  53 = 0
  54 = 0
  55 = 0
  56 = 0
  57 = 0
}

foreach ($row in $af29.Keys) {
  $ws.Cells.Item($row, 32).Value = $af29[$row]
}

# --- row 54 was missing the provider name (nom_upgd); fill it in ------
$ws.Range("C54").Value = "CLINICA MEDICA TURIN SAS"
